# Delete the rows corresponding to stations that were dropped in this update.
# (USACE, 82770) -> row 6
# (USACE, 76220) -> row 17
# (USACE, 76593) -> row 21
# Delete from bottom to top so earlier row numbers stay valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(6).Delete()

# Reset the saved selection/active cell back to A1 (was C1:C33 before the edit).
$ws.Range("A1").Select()
